$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.66533855846194
$ws.Range("C2").Value = 0.2069052540594214
$ws.Range("D2").Value = 0.6670649011318801
$ws.Range("E2").Value = 0.2718899746259495
$ws.Range("G2").Value = 1.706659762225328
$ws.Range("H2").Value = 1.497637982732471
$ws.Range("I2").Value = 1.192978001255106
$ws.Range("J2").Value = 0.1413534793413191
$ws.Range("K2").Value = 0.8615377913285158
$ws.Range("N2").Value = 2.29462689587843
$ws.Range("B3").Value = 0.6255790943371551
$ws.Range("C3").Value = 0.196803071248155
$ws.Range("D3").Value = 0.6572865029390016
$ws.Range("E3").Value = 0.2670595435830023
$ws.Range("G3").Value = 1.700031355704311
$ws.Range("H3").Value = 1.500133733699911
$ws.Range("I3").Value = 1.19479847644493
$ws.Range("J3").Value = 0.1381284282075654
$ws.Range("K3").Value = 0.8127085917121519
$ws.Range("N3").Value = 2.313420468893788
$ws.Range("B4").Value = 0.6015062699190423
$ws.Range("C4").Value = 0.1907191670383668
$ws.Range("D4").Value = 0.6516295603923652
$ws.Range("E4").Value = 0.2642417031785484
$ws.Range("G4").Value = 1.696870084748966
$ws.Range("H4").Value = 1.502289647577683
$ws.Range("I4").Value = 1.19646012682432
$ws.Range("J4").Value = 0.1362298764566034
$ws.Range("K4").Value = 0.7831831601621388
$ws.Range("N4").Value = 2.325651018500324
$ws.Range("B5").Value = 0.5917819473527572
$ws.Range("C5").Value = 0.188269749428116
$ws.Range("D5").Value = 0.6494115604354249
$ws.Range("E5").Value = 0.2631306127431259
$ws.Range("G5").Value = 1.69580998607249
$ws.Range("H5").Value = 1.503324967112064
$ws.Range("I5").Value = 1.197273915937195
$ws.Range("J5").Value = 0.1354766944674139
$ws.Range("K5").Value = 0.7712659949647502
$ws.Range("N5").Value = 2.330808748098832
$ws.Range("B6").Value = 0.590172402243752
$ws.Range("C6").Value = 0.1878648247210606
$ws.Range("D6").Value = 0.6490485325694806
$ws.Range("E6").Value = 0.2629483631246998
$ws.Range("G6").Value = 1.695647727353631
$ws.Range("H6").Value = 1.503506348956577
$ws.Range("I6").Value = 1.197417295408414
$ws.Range("J6").Value = 0.1353528658651584
$ws.Range("K6").Value = 0.7692940872723568
$ws.Range("N6").Value = 2.331675671303465
$ws.Range("B7").Value = 0.6013747777302569
$ws.Range("C7").Value = 0.1906860126768777
$ws.Range("D7").Value = 0.6515992944131597
$ws.Range("E7").Value = 0.2642265680307929
$ws.Range("G7").Value = 1.696854864578782
$ws.Range("H7").Value = 1.502302975548659
$ws.Range("I7").Value = 1.196470548713336
$ws.Range("J7").Value = 0.1362196358518233
$ws.Range("K7").Value = 0.7830219766691755
$ws.Range("N7").Value = 2.32571987440069
$ws.Range("B8").Value = 0.6515590715424366
$ws.Range("C8").Value = 0.2033973025556293
$ws.Range("D8").Value = 0.6636212686698855
$ws.Range("E8").Value = 0.2701936912048666
$ws.Range("G8").Value = 1.704185502637358
$ws.Range("H8").Value = 1.498369092089106
$ws.Range("I8").Value = 1.193492752974073
$ws.Range("J8").Value = 0.1402245133428366
$ws.Range("K8").Value = 0.8446068734805863
$ws.Range("N8").Value = 2.300963379090923
$ws.Range("B9").Value = 0.752665008631368
$ws.Range("C9").Value = 0.2292720640538732
$ws.Range("D9").Value = 0.6899527037633675
$ws.Range("E9").Value = 0.2830723440769489
$ws.Range("G9").Value = 1.725788074391858
$ws.Range("H9").Value = 1.49560433612848
$ws.Range("I9").Value = 1.191974697187632
$ws.Range("J9").Value = 0.1487281130066336
$ws.Range("K9").Value = 0.968997299568656
$ws.Range("N9").Value = 2.257905893508656
$ws.Range("B10").Value = 0.8285987903032037
$ws.Range("C10").Value = 0.2488689719054946
$ws.Range("D10").Value = 0.710985832330465
$ws.Range("E10").Value = 0.2932565610704643
$ws.Range("G10").Value = 1.746094567214499
$ws.Range("H10").Value = 1.496595617483734
$ws.Range("I10").Value = 1.193503826534432
$ws.Range("J10").Value = 0.1553759805843242
$ws.Range("K10").Value = 1.06261364406032
$ws.Range("N10").Value = 2.229625103981213
$ws.Range("B11").Value = 0.8635040115800052
$ws.Range("C11").Value = 0.2579135709413833
$ws.Range("D11").Value = 0.7209225061036193
$ws.Range("E11").Value = 0.2980475650531318
$ws.Range("G11").Value = 1.756302191079698
$ws.Range("H11").Value = 1.497704195228096
$ws.Range("I11").Value = 1.194775994267523
$ws.Range("J11").Value = 0.1584881035638261
$ws.Range("K11").Value = 1.105690314012264
$ws.Range("N11").Value = 2.217489134926502
$ws.Range("B12").Value = 0.8767738586766995
$ws.Range("C12").Value = 0.2613573204668285
$ws.Range("D12").Value = 0.7247383520348478
$ws.Range("E12").Value = 0.2998845976263169
$ws.Range("G12").Value = 1.760307516176198
$ws.Range("H12").Value = 1.498218633026852
$ws.Range("I12").Value = 1.1953408050436
$ws.Range("J12").Value = 0.1596792916164844
$ws.Range("K12").Value = 1.12207295916923
$ws.Range("N12").Value = 2.212998568002021
$ws.Range("B13").Value = 0.87391364924963
$ws.Range("C13").Value = 0.2606148118922249
$ws.Range("D13").Value = 0.723914181578948
$ws.Range("E13").Value = 0.2994879458311388
$ws.Range("G13").Value = 1.759438668045817
$ws.Range("H13").Value = 1.498103629104804
$ws.Range("I13").Value = 1.195215465867953
$ws.Range("J13").Value = 0.1594221826530031
$ws.Range("K13").Value = 1.118541530835131
$ws.Range("N13").Value = 2.213961015805815
$ws.Range("B14").Value = 0.8645946885644662
$ws.Range("C14").Value = 0.2581965139133047
$ws.Range("D14").Value = 0.7212353747461862
$ws.Range("E14").Value = 0.2981982419767846
$ws.Range("G14").Value = 1.756628904965424
$ws.Range("H14").Value = 1.497744621094228
$ws.Range("I14").Value = 1.194820795848869
$ws.Range("J14").Value = 0.1585858486483005
$ws.Range("K14").Value = 1.107036712441669
$ws.Range("N14").Value = 2.21711758644232
$ws.Range("B15").Value = 0.8588933236001139
$ws.Range("C15").Value = 0.2567176816241101
$ws.Range("D15").Value = 0.7196014372095476
$ws.Range("E15").Value = 0.2974112294303453
$ws.Range("G15").Value = 1.754926078157126
$ws.Range("H15").Value = 1.497537045998484
$ws.Range("I15").Value = 1.194589871824178
$ws.Range("J15").Value = 0.1580752244522614
$ws.Range("K15").Value = 1.099998852087225
$ws.Range("N15").Value = 2.219064765681431
$ws.Range("B16").Value = 0.8263249357490565
$ws.Range("C16").Value = 0.2482805062735736
$ws.Range("D16").Value = 0.710343867772707
$ws.Range("E16").Value = 0.2929466436258181
$ws.Range("G16").Value = 1.74544703168317
$ws.Range("H16").Value = 1.496536409101083
$ws.Range("I16").Value = 1.193432302838872
$ws.Range("J16").Value = 0.1551743700487123
$ws.Range("K16").Value = 1.059808341897366
$ws.Range("N16").Value = 2.230432907972521
$ws.Range("B17").Value = 0.8064380458517917
$ws.Range("C17").Value = 0.2431378920442739
$ws.Range("D17").Value = 0.7047590828314014
$ws.Range("E17").Value = 0.2902483005069172
$ws.Range("G17").Value = 1.739880706845724
$ws.Range("H17").Value = 1.496091040785927
$ws.Range("I17").Value = 1.192869955383649
$ws.Range("J17").Value = 0.1534173533466401
$ws.Range("K17").Value = 1.035278272486778
$ws.Range("N17").Value = 2.237593788968489
$ws.Range("B18").Value = 0.7950337669869612
$ws.Range("C18").Value = 0.2401922151268536
$ws.Range("D18").Value = 0.701581553963706
$ws.Range("E18").Value = 0.2887111684634576
$ws.Range("G18").Value = 1.736770388813142
$ws.Range("H18").Value = 1.495896772032012
$ws.Range("I18").Value = 1.192600766727104
$ws.Range("J18").Value = 0.1524150417399568
$ws.Range("K18").Value = 1.021215359999218
$ws.Range("N18").Value = 2.2417811517358
$ws.Range("B19").Value = 0.7911783409929285
$ws.Range("C19").Value = 0.239196955666273
$ws.Range("D19").Value = 0.7005116544609962
$ws.Range("E19").Value = 0.2881932769892899
$ws.Range("G19").Value = 1.735732954096534
$ws.Range("H19").Value = 1.495841624304489
$ws.Range("I19").Value = 1.192518938237953
$ws.Range("J19").Value = 0.1520770966641862
$ws.Range("K19").Value = 1.0164618242637
$ws.Range("N19").Value = 2.243210698720134
$ws.Range("B20").Value = 0.8085515074848502
$ws.Range("C20").Value = 0.2436840669272726
$ws.Range("D20").Value = 0.7053500019838168
$ws.Range("E20").Value = 0.2905340030527981
$ws.Range("G20").Value = 1.740463800956377
$ws.Range("H20").Value = 1.496132044523335
$ws.Range("I20").Value = 1.192924201499252
$ws.Range("J20").Value = 0.1536035337234125
$ws.Range("K20").Value = 1.037884765754995
$ws.Range("N20").Value = 2.236824398162163
$ws.Range("B21").Value = 0.8673304831065991
$ws.Range("C21").Value = 0.2589063169565975
$ws.Range("D21").Value = 0.7220207653467412
$ws.Range("E21").Value = 0.2985764405951556
$ws.Range("G21").Value = 1.757450399741657
$ws.Range("H21").Value = 1.497847501232059
$ws.Range("I21").Value = 1.194934464257926
$ws.Range("J21").Value = 0.1588311553363724
$ws.Range("K21").Value = 1.110414044871646
$ws.Range("N21").Value = 2.216187572047978
$ws.Range("B22").Value = 0.9060489602027815
$ws.Range("C22").Value = 0.2689643015053491
$ws.Range("D22").Value = 0.7332252692539498
$ws.Range("E22").Value = 0.303965457050154
$ws.Range("G22").Value = 1.769367880436903
$ws.Range("H22").Value = 1.499520343418254
$ws.Range("I22").Value = 1.196732549886683
$ws.Range("J22").Value = 0.1623217230893346
$ws.Range("K22").Value = 1.158226717500497
$ws.Range("N22").Value = 2.203312694049096
$ws.Range("B23").Value = 0.8853564989311735
$ws.Range("C23").Value = 0.2635861317243098
$ws.Range("D23").Value = 0.727216909602646
$ws.Range("E23").Value = 0.301077070758403
$ws.Range("G23").Value = 1.762932514728902
$ws.Range("H23").Value = 1.498577010320787
$ws.Range("I23").Value = 1.195728514146019
$ws.Range("J23").Value = 0.1604519543317622
$ws.Range("K23").Value = 1.132670634700673
$ws.Range("N23").Value = 2.21012814386134
$ws.Range("B24").Value = 0.8075959212657722
$ws.Range("C24").Value = 0.2434371073785542
$ws.Range("D24").Value = 0.7050827439074112
$ws.Range("E24").Value = 0.290404792763411
$ws.Range("G24").Value = 1.740199904294911
$ws.Range("H24").Value = 1.496113314286163
$ws.Range("I24").Value = 1.192899508275389
$ws.Range("J24").Value = 0.1535193372192651
$ws.Range("K24").Value = 1.036706246325963
$ws.Range("N24").Value = 2.237172020209016
$ws.Range("B25").Value = 0.725023824914274
$ws.Range("C25").Value = 0.222169863869766
$ws.Range("D25").Value = 0.6825335417140082
$ws.Range("E25").Value = 0.2794618286708896
$ws.Range("G25").Value = 1.719167285125209
$ws.Range("H25").Value = 1.495821869013156
$ws.Range("I25").Value = 1.191921711918049
$ws.Range("J25").Value = 0.1463576620546618
$ws.Range("K25").Value = 0.9349566990196649
$ws.Range("N25").Value = 2.268965670551914
